{"js": "// Update the two-digit \u00f7 one-digit division answer table.\n// The document has a single table; only every 4th row (0, 4, 8, 12, 16)\n// carries the 5 answer cells per row - the rows in between are spacer rows.\nconst table = context.document.body.tables.getItemAt(0);\n\n// Each entry is [rowIndex, columnIndex, newText].\nconst updates = [\n  [0, 0, \"22\u00f75=4, 2\"],\n  [0, 1, \"65\u00f75=13, 0\"],\n  [0, 2, \"64\u00f75=12, 4\"],\n  [0, 3, \"73\u00f75=14, 3\"],\n  [0, 4, \"12\u00f74=3, 0\"],\n  [4, 0, \"68\u00f73=22, 2\"],\n  [4, 1, \"15\u00f78=1, 7\"],\n  [4, 2, \"95\u00f72=47, 1\"],\n  [4, 3, \"49\u00f77=7, 0\"],\n  [4, 4, \"66\u00f75=13, 1\"],\n  [8, 0, \"23\u00f79=2, 5\"],\n  [8, 1, \"56\u00f73=18, 2\"],\n  [8, 2, \"23\u00f76=3, 5\"],\n  [8, 3, \"24\u00f76=4, 0\"],\n  [8, 4, \"78\u00f76=13, 0\"],\n  [12, 0, \"78\u00f76=13, 0\"],\n  [12, 1, \"39\u00f74=9, 3\"],\n  [12, 2, \"40\u00f75=8, 0\"],\n  [12, 3, \"53\u00f77=7, 4\"],\n  [12, 4, \"46\u00f74=11, 2\"],\n  [16, 0, \"60\u00f78=7, 4\"],\n  [16, 1, \"91\u00f76=15, 1\"],\n  [16, 2, \"17\u00f72=8, 1\"],\n  [16, 3, \"92\u00f79=10, 2\"],\n  [16, 4, \"13\u00f74=3, 1\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit / one-digit division answer table.\n# The document has a single table; only every 4th row (rows 1, 5, 9, 13, 17\n# in 1-based COM indexing) carries the 5 answer cells per row - the rows\n# in between are spacer rows.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each entry is row (1-based), column (1-based), new text.\n$updates = @(\n    @(1, 1, '22\u00f75=4, 2'),\n    @(1, 2, '65\u00f75=13, 0'),\n    @(1, 3, '64\u00f75=12, 4'),\n    @(1, 4, '73\u00f75=14, 3'),\n    @(1, 5, '12\u00f74=3, 0'),\n    @(5, 1, '68\u00f73=22, 2'),\n    @(5, 2, '15\u00f78=1, 7'),\n    @(5, 3, '95\u00f72=47, 1'),\n    @(5, 4, '49\u00f77=7, 0'),\n    @(5, 5, '66\u00f75=13, 1'),\n    @(9, 1, '23\u00f79=2, 5'),\n    @(9, 2, '56\u00f73=18, 2'),\n    @(9, 3, '23\u00f76=3, 5'),\n    @(9, 4, '24\u00f76=4, 0'),\n    @(9, 5, '78\u00f76=13, 0'),\n    @(13, 1, '78\u00f76=13, 0'),\n    @(13, 2, '39\u00f74=9, 3'),\n    @(13, 3, '40\u00f75=8, 0'),\n    @(13, 4, '53\u00f77=7, 4'),\n    @(13, 5, '46\u00f74=11, 2'),\n    @(17, 1, '60\u00f78=7, 4'),\n    @(17, 2, '91\u00f76=15, 1'),\n    @(17, 3, '17\u00f72=8, 1'),\n    @(17, 4, '92\u00f79=10, 2'),\n    @(17, 5, '13\u00f74=3, 1'),\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
